$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old header row (row 2, columns A-H)
$ws.Range("A2:H2").ClearContents()

# Re-write the header, now on row 1, columns A-H, plus two new columns I, J
$ws.Range("A1").Value = "文件名稱"
$ws.Range("B1").Value = "分類"
$ws.Range("C1").Value = "型號 / Model name"
$ws.Range("D1").Value = "正常電壓 / Nominal voltage,單位 V"
$ws.Range("E1").Value = "典型/正常電池容量,Wh"
$ws.Range("F1").Value = "典型/正常容量,mAh"
$ws.Range("G1").Value = "額定容量,mA"
$ws.Range("H1").Value = "額定能量,Wh"
$ws.Range("I1").Value = "備註"
$ws.Range("J1").Value = "衝突"

# Update the selection to match the edited workbook
$ws.Range("J2").Select()
